$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: file names, written directly (plain text, no date risk) ---
$names = @(
  'Routine_Care/Nursing for Arterial and Central Venous Lines.pdf',
  'Routine_Care/VTE_Prevention/TED Stocking Sizing.pdf',
  'Routine_Care/Faecal  incontinence skin care.pdf',
  'Breathing(Respiratory)/Equipment/IPPB using an ICU Ventilator.pdf',
  'Drugs/heparin_critical_care_only.pdf',
  'Routine_Care/Invasive Flush Systems.pdf',
  'GI_Liver_and_Transplant/Pancreatic Irrigation.pdf',
  'Neurological/Thiopentone levels.pdf',
  'Infection_and_sepsis/Ebola/Ebola.pdf',
  'GI_Liver_and_Transplant/Nasal bridle.pdf',
  'Breathing(Respiratory)/Equipment/AMBU AScope.pdf',
  'Cardiovascular/EZ-IO Intraosseus Access Device_pub_em.pdf',
  'Routine_Care/Central venous catheter removal.pdf',
  'Routine_Care/Tracheostomy_nursing_care.pdf',
  'Airway/Critical care extubation checklist.pdf',
  'ECLS/Extra Corporeal Carbon Dioxide Removal.pdf',
  'Airway/Tracheostomy_Laryngectomy/Hospital_in-patients_with_a_Tracheostomy.pdf',
  'Trauma and Burns/Mangement of burns.pdf',
  'Drugs/diazepam_diazemuls.pdf',
  'Breathing(Respiratory)/salbutamol and ipratroprium MDI.pdf',
  'End_of_life_care/Reasons to report a death to PF.pdf',
  'Neurological/SOP -  Femoral site care.pdf',
  'Airway/Tracheostomy_Laryngectomy/Tracheostomy change in Critical Care.pdf',
  'Airway/Tracheostomy_Laryngectomy/Tracheostomy suctioning cleaning guideline.pdf',
  'Drugs/sodium_bicarbonate.pdf',
  'Post_op_care/Anticoagulation antiplatelet agents and epidural analgesia.pdf',
  'Post_op_care/Epidural top-up.pdf',
  'Diabetes_and_Glucose/Hyperosmolar Hyperglycaemic State.pdf',
  'Drugs/heparin for Haemofiltration.pdf',
  'Covid-19/SJH/SJH COVID19 ED Intubation Action Card.pdf',
  'Covid-19/SJH/SJH COVID19 ITU Intubation Action Card.pdf',
  'Covid-19/WGH/CoVid intubation checklist WGH.pdf',
  'Airway/Emergency intubation checklist_em_pub.pdf',
  'Drugs/fentanyl.pdf',
  'Airway/Tracheostomy_Laryngectomy/Tracheostomy guideline.pdf',
  'Covid-19/WGH/WGH_CT_Transfer_May.pdf',
  'Cardiovascular/GJNH Acute Heart Failure Referral Form.pdf',
  'Organ_donation/Donation after circulatory death.pdf',
  'Airway/Percutaneous tracheostomy checklist.pdf',
  'Delirium/Managing a Potentially Violent Patient.pdf',
  'Delirium/Risk assessment posi mit.pdf',
  'Infection_and_sepsis/SOP Ultrasound Cleaning.pdf',
  'GI_Liver_and_Transplant/Treatment of constipation.pdf',
  'GI_Liver_and_Transplant/Abdominal pressure measurement.pdf',
  'Airway/Anticipated difficult airway tool.pdf',
  'Drugs/ketamine_in_asthma.pdf',
  'Breathing(Respiratory)/HFNO.pdf',
  'Airway/McGrath Mac.pdf',
  'Airway/Tracheostomy_Laryngectomy/Tracheostomy safety box contents.pdf',
  'Delirium/Drugs Causing Delirium and Agitiation.pdf',
  'Neurological/Sub arachnoid haemorrhage management.pdf',
  'End_of_life_care/Documentation following death.pdf',
  'Drugs/zanamivir.pdf',
  'Routine_Care/bBraun Spaceplus Failure EMERGENCY ACTION CARD_em.pdf',
  'Drugs/insulin.pdf',
  'Breathing(Respiratory)/Equipment/HFNO Set Up.pdf',
  'Breathing(Respiratory)/Inhaled Nitrous Oxide.pdf',
  'Cardiovascular/Steroids for Septic Shock.pdf',
  'Breathing(Respiratory)/Equipment/APRV.pdf',
  'Post_op_care/Epidural Haematoma.pdf',
  'Neurological/SOP for review of Neurosurgical patients in ITU by neurosurgical team.pdf',
  'Breathing(Respiratory)/Equipment/T piece Y piece.pdf',
  'Policies_and_admin/General Critical Care Interaction with HEPMA_pub.pdf',
  'Drugs/midazolam and thiopental levels.pdf',
  'Breathing(Respiratory)/Equipment/HFNO through ventilator.pdf',
  'Routine_Care/VTE_Prevention/Dalteparin_thromboprophylaxis.pdf',
  'Post_op_care/Adult Scoliosis Spinal Surgery Post-Op Care.pdf',
  'Post_op_care/Post op care pharyngo-laryngo-oesphagectomy PLOG.pdf',
  'GI_Liver_and_Transplant/Nasogastric feeding protocol.pdf',
  'Drugs/Antibiotic doses in CVVHD.pdf',
  'Diabetes_and_Glucose/Intravenous Insulin Therapy (not for DKA or HHS).pdf',
  'GI_Liver_and_Transplant/Nasojejunal feeding protocol.pdf',
  'GI_Liver_and_Transplant/Jejunostomy feeding protocol.pdf',
  'ECLS/RIE ECLS Anti Xa Protocol.pdf',
  'Infection_and_sepsis/Winter Infections Stepdown Guidance.pdf',
  'Drugs/vasopressin_sepsis.pdf',
  'Drugs/vasopressin organ donation.pdf',
  'Covid-19/videos/Donning and Doffing Video.pdf',
  'Transfer/ACCP Transfers.pdf',
  'Breathing(Respiratory)/Equipment/Bipap V60.pdf',
  'Breathing(Respiratory)/CPAP.pdf',
  'Breathing(Respiratory)/Equipment/Ventilators Circuits Filters and Closed Suction - Set up and Maintenance.pdf',
  'Infection_and_sepsis/Infection indications for IVIG.pdf',
  'Drugs/piperacillin_tazobactam extended_infusion.pdf',
  'Procedures/CVC Guidance/Securing CVCs.pdf',
  'Covid-19/Covid 19 Death Certification Guideline.pdf',
  'Neurological/Treatment of status epilepticus.pdf',
  'Routine_Care/Video Communication.pdf',
  'Cardiovascular/Cardiogenic Shock.pdf',
  'Drugs/isoprenaline.pdf',
  'Haematology_CAR-T/Haem_ICU_transfer.pdf',
  'Drugs/aminophylline.pdf',
  'Cardiovascular/Management of hypertension within Critical Care.pdf',
  'Haematology_CAR-T/CRS.pdf',
  'Drugs/pancuronium.pdf',
  'Haematology_CAR-T/ICANS.pdf',
  'Drugs/phenytoin.pdf',
  'Drugs/rocuronium.pdf',
  'Drugs/milrinone.pdf',
  'Policies_and_admin/General Critical Care SOP_pub.pdf',
  'Neurological/Management of traumatic brain injury.pdf',
  'Ethics_and_Law/DNACPR policy for Scotland.pdf',
  'Ethics_and_Law/Care at the End of Life (FICM).pdf',
  'End_of_life_care/CMO & NRS Guidance for Doctors completing MCCD.pdf',
  'GI_Liver_and_Transplant/ICU - Upper GI bleeding (Endoscopy guideline).pdf',
  'Neurological/Critical Care MRI Procedure_pub.pdf',
  'Infection_and_sepsis/Trip Out of Unit infection guidance.pdf',
  'Covid-19/COVID 19 ICM guidance basic goals_June_2022.pdf',
  'End_of_life_care/Palliative extubation & withdrawal of invasive ventilatory support nursing checklist.pdf',
  'Organ_donation/Organ Retrieval SOP.pdf',
  'Drugs/clonidine.pdf',
  'Airway/Cook Staged Extubation Set.pdf',
  'Drugs/noradrenaline (central).pdf',
  'Post_op_care/Epidural hypotension.pdf',
  'Breathing(Respiratory)/Equipment/Passy Muir Valve.pdf',
  'GI_Liver_and_Transplant/Confirmation of Nasogastric Tube Position.pdf',
  'Drugs/dexmedetomidine.pdf',
  'Drugs/glyceryl_trinitrate.pdf',
  'GI_Liver_and_Transplant/Fulminant Liver Failure.pdf',
  'Infection_and_sepsis/Antifungal guidance in critical care.pdf',
  'Cardiovascular/Intra Aortic Balloon Pump Guideline_pub.pdf',
  'Cardiovascular/Intra Aortic Balloon Pump Bedside Checks_pub.pdf',
  'Transfer/Transfer Outdoors to Garden Guideline.pdf',
  'Breathing(Respiratory)/ARDS Strategy.pdf',
  'Drugs/dobutamine.pdf',
  'Drugs/adrenaline.pdf',
  'Drugs/hydralazine.pdf',
  'Post_op_care/Major OMFS Free Flap.pdf',
  'Drugs/alfentanil.pdf',
  'Drugs/Alteplase for massive PE.pdf',
  'Drugs/magnesium.pdf',
  'Drugs/ICU - IV drug infusions.pdf',
  'Drugs/neostigmine.pdf',
  'Drugs/vancomycin.pdf',
  'Drugs/labetalol.pdf',
  'Neurological/Intrathecal policy RIE.pdf',
  'Infection_and_sepsis/Initial investigation and management in unidentified Infections.pdf',
  'Drugs/midazolam.pdf',
  'Cardiovascular/Management of Acute Type B Aortic Dissection Guideline.pdf',
  'Drugs/potassium.pdf',
  'Procedures/CVC Guidance/CVC NHL  April 2023.pdf',
  'Drugs/salbutamol.pdf',
  'Drugs/nimodipine.pdf',
  'Drugs/phenobarbitone.pdf',
  'Drugs/nicardipine.pdf',
  'Routine_Care/ICU Eye Care Guideline.pdf',
  'Procedures/Arterial Line insertion for ACCPs.pdf',
  'Drugs/amiodarone.pdf',
  'Drugs/phenylephrine.pdf',
  'Breathing(Respiratory)/Manual Ventilation and MHI.pdf',
  'Drugs/noradrenaline (peripheral).pdf',
  'Neurological/Ventriculitis Guideline.pdf',
  'Drugs/Epoprostenol.pdf',
  'Drugs/morphine.pdf',
  'Cardiovascular/Cardiac Output Monitoring_pub .pdf',
  'Cardiovascular/Pulmonary_Embolism_and_DVT/Catheter directed thrombolysis of iliofemoral DVT alteplase_pub.pdf',
  'Drugs/calcium.pdf',
  'Drugs/dalteparin_thromboprophylaxis.pdf',
  'Drugs/Vancomycin Continuous Infusion Fluid Restricted.pdf',
  'Drugs/atracurium.pdf',
  'Airway/Tracheostomy_Laryngectomy/Decannulation Guidline.pdf',
  'Policies_and_admin/Pet Visitation.pdf',
  'Drugs/valproate.pdf',
  'Procedures/Inadvertent Catheter Placement Guideline.pdf',
  'Drugs/stress ulcer prophylaxis.pdf',
  'Drugs/ketamine_for_status epilepticus.pdf',
  'Drugs/Phosphate.pdf',
  'Drugs/Thiopentone.pdf',
  'Breathing(Respiratory)/Proning Guideline.pdf',
  'GI_Liver_and_Transplant/Prokinetics in ICU.pdf',
  'Drugs/Octreotide.pdf',
  'Procedures/ACCP CVC placement following completion of initial competencies.pdf',
  'Procedures/ACCPs acquiring initial CVC competencies.pdf',
  'Post_op_care/Prevention and treatment of paraplegia after major aortic procedures.pdf',
  'Delirium/Violence and Agression.pdf',
  'Post_op_care/Care of the Transgender Patient.pdf',
  'GI_Liver_and_Transplant/Plasma exchange in Acute Liver Failure.pdf',
  'End_of_life_care/Guideline following Sudden Cardiac Death where death occurs in ICU.pdf',
  'Breathing(Respiratory)/Equipment/NIV through Drager Vent Set up in Critical Care.pdf',
  'Transfer/Transfer Guidelines.pdf',
  'Breathing(Respiratory)/Equipment/NIV through Nihon Kohden  Setup.pdf',
  'Infection_and_sepsis/Influenza in Critical Care.pdf',
  'Policies_and_admin/Anticipated Post op flow surgical patients.pdf',
  'Policies_and_admin/Discharge Home from Critical Care.pdf',
  'Policies_and_admin/Repatriaiton Checklist for Critical Care.pdf'
)

# --- Column B: review dates, written as literal text via TEXT() + paste-values ---
# (direct string assignment would be auto-parsed into a date serial by Excel)
$dates = @(
  '2011-03-07',
  '2011-06-07',
  '2011-07-07',
  '2011-12-07',
  '2014-09-07',
  '2014-10-07',
  '2014-11-07',
  '2014-12-07',
  '2015-03-07',
  '2016-01-07',
  '2016-04-07',
  '2016-08-07',
  '2017-01-07',
  '2017-08-07',
  '2017-10-07',
  '2017-10-07',
  '2017-11-07',
  '2018-05-07',
  '2018-07-07',
  '2019-05-07',
  '2019-05-07',
  '2019-06-07',
  '2019-06-07',
  '2019-06-07',
  '2019-08-07',
  '2019-09-07',
  '2020-01-07',
  '2020-03-07',
  '2020-03-07',
  '2020-03-07',
  '2020-03-07',
  '2020-03-07',
  '2020-03-07',
  '2020-04-07',
  '2020-05-07',
  '2020-07-07',
  '2020-08-07',
  '2020-11-07',
  '2021-02-07',
  '2021-05-07',
  '2021-05-07',
  '2021-05-07',
  '2021-06-07',
  '2021-06-07',
  '2021-06-07',
  '2021-06-07',
  '2021-06-07',
  '2021-06-07',
  '2021-06-07',
  '2021-06-07',
  '2021-06-07',
  '2021-09-07',
  '2021-12-07',
  '2022-01-07',
  '2022-03-07',
  '2022-03-07',
  '2022-04-07',
  '2022-05-07',
  '2022-05-07',
  '2022-06-07',
  '2022-06-07',
  '2022-06-07',
  '2022-07-07',
  '2022-08-07',
  '2022-10-07',
  '2022-11-07',
  '2022-11-07',
  '2022-12-07',
  '2023-01-07',
  '2023-02-07',
  '2023-02-07',
  '2023-04-07',
  '2023-04-07',
  '2023-04-07',
  '2023-05-07',
  '2023-05-07',
  '2023-05-07',
  '2023-06-07',
  '2023-06-07',
  '2023-07-07',
  '2023-07-07',
  '2023-07-07',
  '2023-07-07',
  '2023-07-07',
  '2023-08-07',
  '2023-08-07',
  '2023-09-07',
  '2023-09-07',
  '2023-10-07',
  '2023-10-07',
  '2024-01-07',
  '2024-02-07',
  '2024-02-07',
  '2024-03-07',
  '2024-03-07',
  '2024-03-07',
  '2024-03-07',
  '2024-03-07',
  '2024-04-07',
  '2024-04-07',
  '2024-05-07',
  '2024-05-07',
  '2024-05-07',
  '2024-05-07',
  '2024-05-07',
  '2024-05-07',
  '2024-05-07',
  '2024-05-07',
  '2024-05-07',
  '2024-05-07',
  '2024-05-07',
  '2024-06-02',
  '2024-06-07',
  '2024-06-26',
  '2024-07-07',
  '2024-07-07',
  '2024-07-07',
  '2024-07-07',
  '2024-07-07',
  '2024-07-25',
  '2024-08-07',
  '2024-08-07',
  '2024-08-07',
  '2024-08-15',
  '2024-10-07',
  '2024-10-07',
  '2024-10-24',
  '2024-11-20',
  '2024-11-24',
  '2024-11-24',
  '2024-12-24',
  '2025-01-07',
  '2025-01-07',
  '2025-01-07',
  '2025-02-07',
  '2025-02-07',
  '2025-02-07',
  '2025-03-07',
  '2025-03-07',
  '2025-04-07',
  '2025-04-07',
  '2025-04-07',
  '2025-04-07',
  '2025-05-07',
  '2025-05-07',
  '2025-05-07',
  '2025-05-07',
  '2025-05-07',
  '2025-06-07',
  '2025-06-07',
  '2025-06-07',
  '2025-06-07',
  '2025-06-07',
  '2025-06-07',
  '2025-07-07',
  '2025-07-07',
  '2025-07-07',
  '2025-08-07',
  '2025-08-07',
  '2025-08-07',
  '2025-08-07',
  '2025-10-07',
  '2025-10-07',
  '2025-10-07',
  '2025-11-07',
  '2025-11-07',
  '2025-11-07',
  '2026-01-07',
  '2026-01-07',
  '2026-01-07',
  '2026-01-07',
  '2026-03-07',
  '2026-03-07',
  '2026-03-07',
  '2026-05-07',
  '2026-07-07',
  '2026-11-07',
  '2027-01-07',
  '2027-01-07',
  '2027-02-07',
  '2028-02-07',
  '2028-05-07',
  '2028-07-07',
  '2028-11-07',
  '2030-11-07'
)

$n = $names.Count
for ($i = 0; $i -lt $n; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $names[$i]
  $ws.Cells.Item($row, 26).Formula = '=TEXT("' + $dates[$i] + '","yyyy-mm-dd")'
}

$helper = $ws.Range($ws.Cells.Item(2, 26), $ws.Cells.Item($n + 1, 26))
$target = $ws.Range($ws.Cells.Item(2, 2), $ws.Cells.Item($n + 1, 2))
$helper.Copy()
$target.PasteSpecial(-4163)
$helper.Clear()
$excel.CutCopyMode = 0

Write-Output "Done. A2=$($ws.Range('A2').Value2) B2=$($ws.Range('B2').Value2)"
